$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Update status values (column F) per the audit update
$ws.Range("E3").Value = 2.5
$ws.Range("F3").Value = "Terminé"
$ws.Range("F7").Value = "Attribué"
$ws.Range("F12").Value = "En cours"
$ws.Range("F17").Value = "En cours"

# Recalculate so the SUM formula in E21 reflects the new value
$excel.Calculate()

# Move the active selection to F19 (was F21)
$ws.Range("F19").Select()
